$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 27.10983833333333
$ws.Range("H2").Value = 81.329515
$ws.Range("I2").Value = 0.9284397459331688
$ws.Range("J2").Value = 0.9284397459331687
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 5.419829666666666
$ws.Range("N2").Value = 16.259489
$ws.Range("O2").Value = 0.5710334161275863
$ws.Range("P2").Value = 0.5710334161275863
$ws.Range("Q2").Value = 146.9307060575372
$ws.Range("R2").Value = 1322.376354517835
$ws.Range("S2").Value = 0.5301701197888457
$ws.Range("T2").Value = 0.5301701197888457

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 27.10983833333333
$ws.Range("H3").Value = 81.329515
$ws.Range("I3").Value = 0.9284397459331688
$ws.Range("J3").Value = 0.9284397459331687
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.326362
$ws.Range("N3").Value = 9.979085999999999
$ws.Range("O3").Value = 0.350465600020454
$ws.Range("P3").Value = 0.3504656000204539
$ws.Range("Q3").Value = 90.17713605814332
$ws.Range("R3").Value = 811.5942245232899
$ws.Range("S3").Value = 0.3253861926413059
$ws.Range("T3").Value = 0.3253861926413057

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 27.10983833333333
$ws.Range("H4").Value = 81.329515
$ws.Range("I4").Value = 0.9284397459331688
$ws.Range("J4").Value = 0.9284397459331687
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.7450736666666667
$ws.Range("N4").Value = 2.235221
$ws.Range("O4").Value = 0.07850098385195992
$ws.Range("P4").Value = 0.0785009838519599
$ws.Range("Q4").Value = 20.19882664975722
$ws.Range("R4").Value = 181.789439847815
$ws.Range("S4").Value = 0.07288343350301744
$ws.Range("T4").Value = 0.07288343350301743

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.190813333333333
$ws.Range("H5").Value = 3.57244
$ws.Range("I5").Value = 0.04078218449921273
$ws.Range("J5").Value = 0.04078218449921273
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.419829666666666
$ws.Range("N5").Value = 16.259489
$ws.Range("O5").Value = 0.5710334161275863
$ws.Range("P5").Value = 0.5710334161275863
$ws.Range("Q5").Value = 6.454005431462221
$ws.Range("R5").Value = 58.08604888315999
$ws.Range("S5").Value = 0.02328799013173094
$ws.Range("T5").Value = 0.02328799013173094

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.190813333333333
$ws.Range("H6").Value = 3.57244
$ws.Range("I6").Value = 0.04078218449921273
$ws.Range("J6").Value = 0.04078218449921273
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.326362
$ws.Range("N6").Value = 9.979085999999999
$ws.Range("O6").Value = 0.350465600020454
$ws.Range("P6").Value = 0.3504656000204539
$ws.Range("Q6").Value = 3.961076221093333
$ws.Range("R6").Value = 35.64968598983999
$ws.Range("S6").Value = 0.01429275276066145
$ws.Range("T6").Value = 0.01429275276066144

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.190813333333333
$ws.Range("H7").Value = 3.57244
$ws.Range("I7").Value = 0.04078218449921273
$ws.Range("J7").Value = 0.04078218449921273
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.7450736666666667
$ws.Range("N7").Value = 2.235221
$ws.Range("O7").Value = 0.07850098385195992
$ws.Range("P7").Value = 0.0785009838519599
$ws.Range("Q7").Value = 0.8872436565822223
$ws.Range("R7").Value = 7.98519290924
$ws.Range("S7").Value = 0.003201441606820349
$ws.Range("T7").Value = 0.003201441606820348

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8986996666666666
$ws.Range("H8").Value = 2.696099
$ws.Range("I8").Value = 0.03077806956761847
$ws.Range("J8").Value = 0.03077806956761847
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 5.419829666666666
$ws.Range("N8").Value = 16.259489
$ws.Range("O8").Value = 0.5710334161275863
$ws.Range("P8").Value = 0.5710334161275863
$ws.Range("Q8").Value = 4.870799114823443
$ws.Range("R8").Value = 43.83719203341099
$ws.Range("S8").Value = 0.01757530620700968
$ws.Range("T8").Value = 0.01757530620700968

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8986996666666666
$ws.Range("H9").Value = 2.696099
$ws.Range("I9").Value = 0.03077806956761847
$ws.Range("J9").Value = 0.03077806956761847
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.326362
$ws.Range("N9").Value = 9.979085999999999
$ws.Range("O9").Value = 0.350465600020454
$ws.Range("P9").Value = 0.3504656000204539
$ws.Range("Q9").Value = 2.989400420612666
$ws.Range("R9").Value = 26.904603785514
$ws.Range("S9").Value = 0.01078665461848668
$ws.Range("T9").Value = 0.01078665461848668

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8986996666666666
$ws.Range("H10").Value = 2.696099
$ws.Range("I10").Value = 0.03077806956761847
$ws.Range("J10").Value = 0.03077806956761847
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.7450736666666667
$ws.Range("N10").Value = 2.235221
$ws.Range("O10").Value = 0.07850098385195992
$ws.Range("P10").Value = 0.0785009838519599
$ws.Range("Q10").Value = 0.6695974558754444
$ws.Range("R10").Value = 6.026377102879
$ws.Range("S10").Value = 0.002416108742122117
$ws.Range("T10").Value = 0.002416108742122117

